# Add 2022-Q3 data
#
# 1) Insert a new worksheet named "2022-Q3" right before the existing
#    "2022-Q2" sheet, cloned from "2022-Q2" so it keeps the same column
#    layout/styles, then overwrite its cells with the 2022-Q3 figures and
#    trim the now-unused trailing rows.
# 2) Update the "总计" (summary) sheet: write a new top data row for
#    2022-Q3 and shift every other quarter's figures down by one row,
#    appending the final (2020-Q4) row at the bottom.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: clone "2022-Q2" -> new sheet placed immediately before it, then
# rename to "2022-Q3".
# ---------------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")
$q2.Copy($q2)
$q3 = $wb.Worksheets.Item("2022-Q2 (2)")
$q3.Name = "2022-Q3"

# The clone inherited 12 data rows from "2022-Q2"; the 2022-Q3 table only
# needs 7 fund rows (rows 2-8). Drop the leftover rows 9-12.
$q3.Rows.Item(9).Resize(4).Delete()

# B..G hold text (fund code/name/scale/position figures are stored as
# strings in this workbook, not numbers) - force a text format before
# writing so values like "009007" or "0.80" keep their exact text form.
$q3.Range("B2:G8").NumberFormat = "@"

# Fund rows: index, code, name, scale, stock position, position pct,
# market value (亿元), position rank.
$q3.Range("A2").Value = 0
$q3.Range("B2").Value = "009007"
$q3.Range("C2").Value = "兴全沪港深两年持有期混合"
$q3.Range("D2").Value = "16.31"
$q3.Range("E2").Value = "92.09"
$q3.Range("F2").Value = "3.40"
$q3.Range("G2").Value = "0.5545"
$q3.Range("H2").Value = 10

$q3.Range("A3").Value = 1
$q3.Range("B3").Value = "009017"
$q3.Range("C3").Value = "银华港股通精选股票A"
$q3.Range("D3").Value = "0.80"
$q3.Range("E3").Value = "80.26"
$q3.Range("F3").Value = "4.12"
$q3.Range("G3").Value = "0.0330"
$q3.Range("H3").Value = 7

$q3.Range("A4").Value = 2
$q3.Range("B4").Value = "501303"
$q3.Range("C4").Value = "广发恒生中型股指数（LOF）A"
$q3.Range("D4").Value = "0.21"
$q3.Range("E4").Value = "89.12"
$q3.Range("F4").Value = "1.47"
$q3.Range("G4").Value = "0.0031"
$q3.Range("H4").Value = 9

$q3.Range("A5").Value = 3
$q3.Range("B5").Value = "004996"
$q3.Range("C5").Value = "广发恒生中型股指数（LOF）C"
$q3.Range("D5").Value = "0.09"
$q3.Range("E5").Value = "89.12"
$q3.Range("F5").Value = "1.47"
$q3.Range("G5").Value = "0.0013"
$q3.Range("H5").Value = 9

$q3.Range("A6").Value = 4
$q3.Range("B6").Value = "160922"
$q3.Range("C6").Value = "大成恒生综合中小型股指数（QDII-LOF）A"
$q3.Range("D6").Value = "0.09"
$q3.Range("E6").Value = "86.62"
$q3.Range("F6").Value = "1.05"
$q3.Range("G6").Value = "0.0009"
$q3.Range("H6").Value = 9

$q3.Range("A7").Value = 5
$q3.Range("B7").Value = "014052"
$q3.Range("C7").Value = "银华港股通精选股票C"
$q3.Range("D7").Value = "0.02"
$q3.Range("E7").Value = "80.26"
$q3.Range("F7").Value = "4.12"
$q3.Range("G7").Value = "0.0008"
$q3.Range("H7").Value = 7

$q3.Range("A8").Value = 6
$q3.Range("B8").Value = "008972"
$q3.Range("C8").Value = "大成恒生综合中小型股指数C"
$q3.Range("D8").Value = "0.02"
$q3.Range("E8").Value = "86.62"
$q3.Range("F8").Value = "1.05"
$q3.Range("G8").Value = "0.0002"
$q3.Range("H8").Value = 9

# Re-apply the sheet's pristine default cell style to B2:G8 so that
# forcing the "@" text format above doesn't leave a stray style index on
# these cells (they should stay on the sheet's default/unstyled xf, same
# as every other data cell outside column A/row 1).
$q2.Range("H2").Copy()
$q3.Range("B2:G8").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# Step 2: update the "总计" summary sheet.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

# Append a new row 9 for 2020-Q4 (was row 8), copying row 8's formatting
# (border/bold index style) onto the new A9 cell first.
$summary.Range("A9").Value = 7
$summary.Range("A8").Copy()
$summary.Range("A9").PasteSpecial(-4122)
$summary.Range("B9").Value = "2020-Q4"
$summary.Range("C9").Value = 6
$summary.Range("D9").Value = 4.21

# Shift the remaining quarters' figures down by one row (bottom-up so we
# never overwrite a source row before it has been copied).
$summary.Range("B8").Value = "2021-Q1"
$summary.Range("C8").Value = 8
$summary.Range("D8").Value = 3.06

$summary.Range("B7").Value = "2021-Q2"
$summary.Range("C7").Value = 11
$summary.Range("D7").Value = 6.63

$summary.Range("B6").Value = "2021-Q3"
$summary.Range("C6").Value = 25
$summary.Range("D6").Value = 13.49

$summary.Range("B5").Value = "2021-Q4"
$summary.Range("C5").Value = 23
$summary.Range("D5").Value = 13.01

$summary.Range("B4").Value = "2022-Q1"
$summary.Range("C4").Value = 15
$summary.Range("D4").Value = 18.12

$summary.Range("B3").Value = "2022-Q2"
$summary.Range("C3").Value = 11
$summary.Range("D3").Value = 1.98

# New top row: 2022-Q3.
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 7
$summary.Range("D2").Value = 0.59

# Restore the original active sheet/selection (sheet-copy operations
# above shift focus onto the newly created sheet).
$summary.Activate()
$null = $summary.Range("A1").Select()
